$d = $word.ActiveDocument

# 1. Replace "17 mins." with "Done" in the third tutorial hyperlink's trailing text.
$d.Content.Find.Execute("17 mins.", $false, $false, $false, $false, $false, $true, 1, $false, "Done", 2) | Out-Null

# 2. Insert a new "Working" tutorial link paragraph right after the "y5eLukU5ur8" one.
$pYur8 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*y5eLukU5ur8*") {
        $pYur8 = $d.Paragraphs($i)
        break
    }
}
$pYur8.Range.InsertParagraphAfter() | Out-Null
$newIdx = $pYur8.Index + 1
$newPara = $d.Paragraphs($newIdx)
$newRange = $newPara.Range
$newRange.Collapse(1)
$d.Hyperlinks.Add($newRange, "https://www.youtube.com/watch?v=RnQvAZHArL0", [Type]::Missing, [Type]::Missing, "https://www.youtube.com/watch?v=RnQvAZHArL0") | Out-Null
$afterLinkRange = $d.Paragraphs($newIdx).Range
$afterLinkRange.Collapse(0)
$afterLinkRange.InsertBefore("  (Working)  ")

# 3. Insert a new "New Hooks" bullet right after "Adding RAG to the System".
$pRag = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "Adding RAG to the System*") {
        $pRag = $d.Paragraphs($i)
        break
    }
}
$pRag.Range.InsertParagraphAfter() | Out-Null
$hooksIdx = $pRag.Index + 1
$hooksRange = $d.Paragraphs($hooksIdx).Range
$hooksRange.Collapse(1)
$hooksRange.InsertBefore("New Hooks like useRef, useEffect")

# 4. Move the _GoBack bookmark from the title paragraph down onto the (still) empty
#    paragraph that now follows the "New Hooks" bullet.
$emptyPara = $null
for ($i = $hooksIdx + 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq "" -or $d.Paragraphs($i).Range.Text -eq [char]13) {
        $emptyPara = $d.Paragraphs($i)
        break
    }
}
$d.Bookmarks.Add("_GoBack", $emptyPara.Range) | Out-Null

Write-Host "Done editing"
